# "make glock/wilson heavier, p320 mag funnel lighter"
# Adds two new barrel rows (Wilson Combat EDC X9) after the existing data,
# separated by one blank spacer row (matching the existing layout pattern
# used between weapon-family blocks: rows 8/11/19/22/25/29/32), and extends
# the three shared formula columns (N, S, T) down to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38: Wilson Combat EDC X9 102mm 9x19 ---
$ws.Range("A38").Value = "wilson_combat_edc_x9_102mm_9x19_barrel"
$ws.Range("B38").Value = "Wilson Combat EDC X9 102mm 9x19"
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 0.16
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = -78
$ws.Range("M38").Value = 0
$ws.Range("P38").Value = 0.06
$ws.Range("Q38").Value = 4.0157499999999997

# --- Row 39: Wilson Combat EDC X9 127mm Threaded 9x19 ---
$ws.Range("A39").Value = "wilson_combat_edc_x9_127mm_threaded_9x19_barrel"
$ws.Range("B39").Value = "Wilson Combat EDC X9 127mm Threaded 9x19"
$ws.Range("C39").Value = -1
$ws.Range("D39").Value = 0.18
$ws.Range("E39").Value = 2
$ws.Range("F39").Value = 2
$ws.Range("H39").Value = 0.1
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = -44
$ws.Range("M39").Value = 750
$ws.Range("P39").Value = 0.06
$ws.Range("Q39").Value = 5

# --- Extend the three shared formula columns down through the new rows,
#     starting a fresh shared-formula group at the row-37 spacer (row 37
#     itself stays otherwise blank, same as the other spacer rows). ---
$ws.Range("N37:N39").Formula = "=C37-D37*20-E37*0.8-F37*0.6-H37*7.5+I37*15+J37/300"
$ws.Range("S37:S39").Formula = "=ROUND(Q37*0.023+P37+R37, 2)"
$ws.Range("T37:T39").Formula = "=(Q37-5)*0.09/11"

# Match the author's final selection position recorded in the saved file.
$ws.Range("D40").Select() | Out-Null
